# Update "想去人数" (interested count) values that changed between scrapes.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 1097
$wsExpo.Range("F4").Value = 1745
$wsExpo.Range("F6").Value = 62
$wsExpo.Range("F7").Value = 199

# Sheet "全部类型" (All types) mirrors the same events
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 1097
$wsAll.Range("F4").Value = 1745
$wsAll.Range("F7").Value = 62
$wsAll.Range("F8").Value = 199
